# Update cryptocurrency price/volume figures (and fix a TRON/WrappedEther row-order
# swap in rows 12-13) per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are forced to Text before/after the write so that
# price strings which happen to look numeric (e.g. "1.001") are not silently
# reinterpreted by Excel as numbers - matching the original inline-string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.849.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.25"
$ws.Range("D3").ClearFormats()
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5042"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3639"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07157"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8913"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.98"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07511"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.221"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008502"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.18"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.901.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.019"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.095.60"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.400"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.782"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.87"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.065"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.09"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.685"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.653"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09160"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05131"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7463"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.977"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.150"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.192"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.557"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01994"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5558"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.069"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.568"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.529"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1468"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4682"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.05"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.553"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.71"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.42%  "
